$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide column B (was merged with column A's width spec; now split with B hidden)
$ws.Columns.Item(2).Hidden = $true

# Update selection (was F10, now C14) and implicitly clear the scrolled topLeftCell
$ws.Range("C14").Select() | Out-Null

# Row 5 - Carry lookahead Adder
$ws.Range("C5").Value = 16560.400000000001
$ws.Range("F5").Value = 2939.6

# Row 7 - Carry Skip Adder
$ws.Range("C7").Value = 16609.099999999999
$ws.Range("F7").Value = 2890.9

# Row 8 - Carry Increment Adder
$ws.Range("C8").Value = 17020.3
$ws.Range("F8").Value = 2479.6999999999998

# Row 9 - Carry bypass Adder
$ws.Range("C9").Value = 15585.5
$ws.Range("F9").Value = 3914.5

# Row 10 - Carry Select Adder
$ws.Range("C10").Value = 17399.400000000001
$ws.Range("F10").Value = 2100.6
$ws.Range("G10").Value = 69
$ws.Range("H10").Value = 239
$ws.Range("I10").Value = 11.118458
$ws.Range("J10").Value = 4.3387460000000004
$ws.Range("K10").Value = 25.687866
